$wb = $excel.ActiveWorkbook

# Update vessels sheet values (column H: max_time_offshore)
$ws = $wb.Worksheets.Item("vessels")
$ws.Range("H2").Value = 120
$ws.Range("H3").Value = 120
$ws.Range("H4").Value = 240

# Select H6 on vessels and make it the active sheet/tab
$ws.Activate()
$ws.Range("H6").Select()
